$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date column (C) for rows 2-8 from 45204 to 45207
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 3).Value = 45207
}
